$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.73547625541687
$ws.Range("B1").Value = 4.038856029510498
$ws.Range("C1").Value = 7.747185230255127
$ws.Range("D1").Value = 7.915988922119141
$ws.Range("E1").Value = 5.81642484664917
